# Insert a new data row above row 282 (shifts existing rows 282:353 down to 283:354)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(282).Insert()

# Populate the newly inserted row 282 with the new record.
$ws.Range("A282").Value = 8
$ws.Range("B282").Value = "Terminal La Palmera de La Serena"
$ws.Range("C282").Value = "Coquimbo"
$ws.Range("D282").Value = 44754
$ws.Range("E282").Value = 4
$ws.Range("F282").Value = 100112032
$ws.Range("G282").Value = "Zapallo italiano"
$ws.Range("H282").Value = "Sin especificar"
$ws.Range("I282").Value = "Primera"
$ws.Range("J282").Value = 500
$ws.Range("K282").Value = 9000
$ws.Range("L282").Value = 10000
$ws.Range("M282").Value = 9500
$ws.Range("N282").Value = "$/caja 50 unidades"
$ws.Range("O282").Value = "Región de Arica y Parinacota"
$ws.Range("P282").Value = 190
$ws.Range("Q282").Value = 50
$ws.Range("R282").Value = "Hortaliza"
